$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("mars")
$ws2 = $wb.Worksheets.Item("april")

# --- mars sheet: add new expense row (row 7) ---
# Force text storage so date-like / numeric-looking strings aren't
# auto-converted to a date serial / number by Excel, then clear the
# temporary number format so no stray cell style is left behind.
$ws1.Range("A7:E7").NumberFormat = "@"
$ws1.Range("A7").Value = "Transportation"
$ws1.Range("B7").Value = "dwas"
$ws1.Range("C7").Value = "2023-03-22"
$ws1.Range("D7").Value = "123.0"
$ws1.Range("E7").Value = "Checkings"
$ws1.Range("A7:E7").ClearFormats()

# --- mars sheet: add monthly total row (row 8) ---
$ws1.Range("A8").Value = "Monthly total: "
$ws1.Range("B8").Value = 5994.0

# --- april sheet: add monthly total row (row 3) ---
$ws2.Range("A3").Value = "Monthly total: "
$ws2.Range("B3").Value = 453.0
